$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: update D, F, G, H
$ws.Cells.Item(38, 4).Value = 2826
$ws.Cells.Item(38, 6).Value = "10.6.21.11"
$ws.Cells.Item(38, 7).Value = 257
$ws.Cells.Item(38, 8).Value = 45519.95326806713

# Row 39: update C, D, E, F, G, H
$ws.Cells.Item(39, 3).Value = 2827
$ws.Cells.Item(39, 4).Value = 3083
$ws.Cells.Item(39, 5).Value = "10.6.21.12"
$ws.Cells.Item(39, 6).Value = "10.6.22.12"
$ws.Cells.Item(39, 7).Value = 257
$ws.Cells.Item(39, 8).Value = 45519.95669515046

# Row 40: update C, D, E, F, G, H
$ws.Cells.Item(40, 3).Value = 3084
$ws.Cells.Item(40, 4).Value = 3340
$ws.Cells.Item(40, 5).Value = "10.6.22.13"
$ws.Cells.Item(40, 6).Value = "10.6.23.13"
$ws.Cells.Item(40, 7).Value = 257
$ws.Cells.Item(40, 8).Value = 45519.95785613426

# Row 41: update C, D, E, F, G, H
$ws.Cells.Item(41, 3).Value = 3341
$ws.Cells.Item(41, 4).Value = 3579
$ws.Cells.Item(41, 5).Value = "10.6.23.14"
$ws.Cells.Item(41, 6).Value = "10.6.23.252"
$ws.Cells.Item(41, 7).Value = 239
$ws.Cells.Item(41, 8).Value = 45519.96153725695

# Row 42: brand new row, same format/style as the block above (row 41's style)
$ws.Range("A41:H41").Copy()
$ws.Range("A42").PasteSpecial(-4122)

$ws.Cells.Item(42, 1).Value = "Домик для Мышки Норушки"
$ws.Cells.Item(42, 2).Value = "без трубы"
$ws.Cells.Item(42, 3).Value = 3580
$ws.Cells.Item(42, 4).Value = 3585
$ws.Cells.Item(42, 5).Value = "10.6.23.253"
$ws.Cells.Item(42, 6).Value = "10.6.24.2"
$ws.Cells.Item(42, 7).Value = 6
$ws.Cells.Item(42, 8).Value = 45519.96697361953
